# Visualizacion de cierre de caja por rol
# Inserts a new "Apertura caja restaurante" menu row into the RolMenu sheet
# right above the existing "Hotel" block (new row 26), shifting the
# following rows down by one, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26 (rows 26-43 shift down to 27-44).
# Excel's default "insert" behaviour copies formatting from the row above
# (row 25), which already carries the s=4 / s=6 / s=5 style pattern we need.
$ws.Rows.Item(26).Insert()

# Populate the new row with the "Apertura caja restaurante" menu entry.
$ws.Range("A26").Value = "Apertura caja restaurante"
$ws.Range("B26").Value = "apertura_caja_restaurante"
$ws.Range("C26").Value = "Restaurante"
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = "lock_open"
$ws.Range("F26").Value = "ADMIN_GERENTE"

# Move the active selection to F27 (was E24 before the edit).
$ws.Range("F27").Select()
